# regenerate orders with updates distance/sizes
#
# The experiment order workbook encodes trial conditions as composite
# labels like "Face05_D80_S20" (Face id / viewing Distance / Size),
# plus matching stimulus filenames ("Face05_D80_S20_l.png" / "_r.png"),
# plus standalone "Distance" (D64/D80/D51) and "Size" (S20/S25/S30)
# columns. This commit regenerates the run with updated distance/size
# values:
#   D64 -> D69
#   D80 -> D86
#   D51 -> D55
#   S30 -> S31
# (S20 and S25 are unchanged.)
#
# Every occurrence of these tokens appears consistently as a substring
# of the Condition / Filename_Left / Filename_Right / Distance / Size
# cell text throughout the sheet, so a straightforward find/replace
# over the whole used range performs the same substitution Excel would
# make if you ran Find & Replace four times.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPart = 2  -> match substrings within a cell's text (not whole-cell only)
# xlByRows = 1
# MatchCase:$true so we don't accidentally touch unrelated lowercase text
$LookAt = 2
$SearchOrder = 1
$MatchCase = $true

# Order matters only in that none of the target tokens are substrings of
# one another, and none of the replacement tokens (D69/D86/D55/S31)
# already occur anywhere in the sheet, so each pass is independent.
$ws.Cells.Replace("D64", "D69", $LookAt, $SearchOrder, $MatchCase)
$ws.Cells.Replace("D80", "D86", $LookAt, $SearchOrder, $MatchCase)
$ws.Cells.Replace("D51", "D55", $LookAt, $SearchOrder, $MatchCase)
$ws.Cells.Replace("S30", "S31", $LookAt, $SearchOrder, $MatchCase)
